# Update the Address column (column C) with the new, shortened address text.
# Row numbers below are worksheet rows (row 1 is the header "Index"/"Name"/"Address").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$addresses = @{
    2  = " Coimbatore, Tamilnadu India"
    3  = " Mahabalipuram, New Delhi - 110074"
    4  = " POLLACHI, COIMBATORE – 642004."
    5  = "Vadodara, Gujarat ,India"
    6  = " New Delhi, Delhi 110068"
    7  = " Thondamuthur, Tamil Nadu 641109"
    8  = "Karjat ,Maharashtra ,India"
    9  = "Mumbai, Maharashtra, India"
    10 = " Amritsar, Punjab 143001"
    11 = "New Delhi, Delhi, India"
    12 = "Palakkad, Kerala,India"
    13 = " Gurgaon - 122001, India"
    14 = "Mulshi, Maharashtra 412115, India"
    15 = "Faridabad, Haryana 121101, India"
    16 = " Bangalore,India -560082"
    17 = "SECTOR 46, GURGAON,India"
    18 = " Noida, Uttar Pradesh, India -201304"
    19 = "Bengaluru, Karnataka, India"
    20 = "New Delhi - 110074"
    21 = "New Delhi- 110074"
    22 = " Gurugram, Haryana, India 122022"
    23 = "New Delhi, India"
    24 = " Gurugram Haryana - 122011`n"
    25 = "Mumbai, Maharashtra, India"
    26 = " Ernakulam, Kerala 682024"
    27 = "Surat India - 395007"
    28 = " Nashik, Maharashtra 422213"
    29 = "Pune, Maharashtra, India"
    30 = "Jaipur, Rajasthan 302026"
    31 = "Powai, Mumbai, Maharashtra 400076, India"
    32 = " Bangalore, Karnataka 560105"
    33 = "Dehradun, Uttarakhand 248001, India"
    34 = " Goa, INDIA 403524"
    35 = "Ghittorni, Delhi, India"
    36 = " Ludhiana, Punjab, India."
    37 = " Rishikesh  India -249307"
    38 = " Kolkata, West Bengal"
    39 = "Kerala, India 680681"
    40 = " Ludhiana, Punjab, India."
}

foreach ($row in $addresses.Keys) {
    $ws.Cells.Item($row, 3).Value = $addresses[$row]
}

# Row 24's address now wraps onto two lines in the source sheet; Excel grew the
# row height and flagged the cell to wrap when the author resized it by hand.
$ws.Cells.Item(24, 3).WrapText = $true
$ws.Rows.Item(24).RowHeight = 28.8

# The author had scrolled down and left the selection on the last data row
# before saving.
$ws.Range("C40").Select()
